$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the two richly-formatted cells (Times New Roman / text number format)
# so they fall back to the sheet's plain default formatting once rewritten. ---
$ws.Range("A2").Clear() | Out-Null
$ws.Range("G2").Clear() | Out-Null
$ws.Range("H2:I2").Clear() | Out-Null

# --- Drop the now-unwanted "nombre_satellites" column (K) entirely. ---
$ws.Columns.Item(11).Delete() | Out-Null

# --- Rename the sheet. ---
$ws.Name = "canteens_good-2"

# --- Normalize formatting on the used range back to the workbook's plain style
# (General number format, Arial 10, no wrap) now that column K is gone. ---
$rng = $ws.Range("A1:J2")
$rng.NumberFormat = "General"
$rng.Font.Name = "Arial"
$rng.Font.Size = 10
$rng.WrapText = $false

# --- Re-enter row 2 data; "siret" becomes a real number instead of text. ---
$ws.Range("A2").Value2 = 21340172201787
$ws.Range("B2").Value2 = "A excel canteen"
$ws.Range("D2").Value2 = 700
$ws.Range("E2").Value2 = 14000
$ws.Range("F2").Value2 = "Cliniques,Hôpitaux"
$ws.Range("G2").Value2 = "Restaurant avec cuisine sur place"
$ws.Range("H2").Value2 = "Concédée"
$ws.Range("I2").Value2 = "Public"

# --- Column widths (converted from OOXML character-width units to the
# Excel ColumnWidth property, which is offset by 5/6 from the stored value). ---
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 13.333333333333332
$ws.Columns.Item(3).ColumnWidth = 15.333333333333334
$ws.Columns.Item(4).ColumnWidth = 16.166666666666668
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666
$ws.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws.Columns.Item(7).ColumnWidth = 27.666666666666668
$ws.Columns.Item(8).ColumnWidth = 10.833333333333332
$ws.Columns.Item(9).ColumnWidth = 17.333333333333336
$ws.Columns.Item(10).ColumnWidth = 21.833333333333336

# --- View settings: zoom out a bit and reset the selection to A1. ---
$excel.ActiveWindow.Zoom = 65
$ws.Range("A1").Select() | Out-Null

# --- Iterative-calc delta (best effort; matches the workbook's calc settings). ---
$excel.Iteration = $false
$excel.MaxChange = 0.001
$excel.MaxIterations = 100
